$d = $word.ActiveDocument
$s = $d.Shapes(1)
Write-Output $s.RelativeHorizontalPosition
Write-Output $s.RelativeVerticalPosition
$s.Left = 2.578503937007874
$s.Top = 1.0
$s.Width = 418.62858267716535
Write-Output "done"
